# Apply cryptos list price/volume updates (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.948.47"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.624.86"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.65"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.250"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.26"
$ws.Range("E10").Value = "  -6.41%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.850.01"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.625.22"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.928.38"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0737"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.23"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.39"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.58"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.134"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.48"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.72"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.18"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.13"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.10"
$ws.Range("E33").Value = "  -5.58%  "
$ws.Range("E34").Value = "  -2.77%  "
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.125.14"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.846"
$ws.Range("E37").Value = "  -6.29%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.518"
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.83"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.767"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.759.83"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.17"
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "54.52"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0528"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  -3.37%  "
